$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 ---------------------------------------------------------------
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 10
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = "2023-11-22 | 14:14 - препарат принят (таблетки)"
$ws.Range("E10").Value = "22/11/2023 14:42:03"

# --- Row 11 -----------------------------------------------------------------
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 11
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = "2023-11-22 | 21:21 - препарат принят ( 1234`n)"
$ws.Range("E11").Value = "22/11/2023 21:00:13"

# Match the existing "id" column formatting (bold, centered, bordered cell
# style) used by A2:A9 for the two freshly added id cells.
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A10:A11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Newline inside D11 makes the engine auto-grow the row height; reset it
# back to the sheet's standard (auto) height so no explicit ht/customHeight
# survives on the exported row, matching un-touched rows elsewhere.
$ws.Rows.Item(10).AutoFit() | Out-Null
$ws.Rows.Item(11).AutoFit() | Out-Null
